$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 52; this shifts existing rows 52-80 down to 53-81
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with the new weekly record
$ws.Range("A52").Value = 8
$ws.Range("B52").Value = "Terminal La Palmera de La Serena"
$ws.Range("C52").Value = "Coquimbo"
$ws.Range("D52").Value = 44455
$ws.Range("E52").Value = 4
$ws.Range("F52").Value = 100112040
$ws.Range("G52").Value = "Cilantro"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 3600
$ws.Range("K52").Value = 2000
$ws.Range("L52").Value = 2500
$ws.Range("M52").Value = 2250
$ws.Range("N52").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O52").Value = "Provincia del Elquí"
$ws.Range("P52").Value = 1500
$ws.Range("Q52").Value = 1.5
$ws.Range("R52").Value = "Hortaliza"

# Match the date display format used by the other rows in column D
$ws.Range("D52").NumberFormat = $ws.Range("D53").NumberFormat
